# Change the HEX suffix back to # in the "stmt"/"prog_line" attribute names
# (column headers and query strings), per commit message: "changed the HEX back to #"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (Table1 column headers) ---
$ws.Range("J1").Value = "stmt#"
$ws.Range("K1").Value = "prog_line#"

# --- Query strings that reference stmtHEX / prog_lineHEX ---
$ws.Range("M14").Value = "stmt s; Select s with s.stmt# = 17"
$ws.Range("M15").Value = "stmt s; Select s with s.stmt# = 99999"
$ws.Range("M16").Value = "stmt s; Select BOOLEAN with s.stmt# = 23"
$ws.Range("M17").Value = "stmt s; Select BOOLEAN with s.stmt# = 999999"
$ws.Range("M18").Value = "assign a; Select a with a.stmt# = 12"
$ws.Range("M19").Value = "assign a; Select a with a.stmt# = 27"
$ws.Range("M20").Value = "assign a; Select BOOLEAN with a.stmt# = 30"
$ws.Range("M21").Value = "assign a; Select BOOLEAN with a.stmt# = 19"
$ws.Range("M22").Value = "while w; Select w with w.stmt# = 19"
$ws.Range("M23").Value = "while w; Select w with w.stmt# = 34"
$ws.Range("M24").Value = "while w; Select BOOLEAN with w.stmt# = 8"
$ws.Range("M25").Value = "while w; Select BOOLEAN with w.stmt# = 25"
$ws.Range("M26").Value = "if i; Select i with i.stmt# = 38"
$ws.Range("M27").Value = "if i; Select i with i.stmt# = 1"
$ws.Range("M28").Value = "if i; Select BOOLEAN with i.stmt# = 38"
$ws.Range("M29").Value = "if i; Select BOOLEAN with i.stmt# = 15"
$ws.Range("M30").Value = "prog_line pl; Select pl with pl.prog_line# = 23"
$ws.Range("M31").Value = "prog_line pl; Select pl with pl.prog_line# = 15081992"
$ws.Range("M32").Value = "prog_line pl; Select BOLLEAN with pl.prog_line# = 1"
$ws.Range("M33").Value = "prog_line pl; Select BOOLEAN with pl.prog_line# = 15081992"

# --- Update window selection to match the author's final cursor position ---
$ws.Range("M45").Select()

Write-Output "done"
